# "chang waitlist logic for active users"
#
# Re-number the waitlist: drop the old StudentID numbers in favor of small
# sequential ids, rename several waitlisted students, and append two more
# students (Reed, Miller) who were waiting behind the existing group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper used below to push a literal (non date-parsed) text value into a
# cell: stage it in a scratch cell formatted as Text, copy it, and paste
# only the *value* into the destination so the destination keeps the
# source's text type instead of Excel re-interpreting "MM/DD/YYYY" as a
# real date.
$scratch = $ws.Range("A200")
function Set-TextValue($rangeAddress, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)  # xlPasteValues
}

# --- Row 2 (Isaac): StudentID becomes a small sequential id ---
$ws.Range("B2").Value = 1

# --- Row 3 (Bea): StudentID becomes a small sequential id ---
$ws.Range("B3").Value = 2

# --- Row 4: renamed from "Rice" to "Lisa" ---
$ws.Range("A4").Value = "Lisa"
Set-TextValue "E4" "09/18/2021"

# --- Row 5: renamed to "Jeff" ---
$ws.Range("A5").Value = "Jeff"
Set-TextValue "E5" "10/18/2021"

# --- Row 6: renamed to "Alice" ---
$ws.Range("A6").Value = "Alice"
Set-TextValue "E6" "11/18/2021"

# --- New row 7: Reed ---
$ws.Range("A7").Value = "Reed"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
Set-TextValue "E7" "12/18/2021"

# --- New row 8: Miller ---
$ws.Range("A8").Value = "Miller"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
Set-TextValue "E8" "01/18/2022"

# Remove the scratch cell so it leaves no trace in the used range.
$ws.Rows.Item(200).Delete()

# Match the saved file's cursor position.
$ws.Range("H9").Select()
